# Resultats_REM-G13006.xlsx fix
# "Arreglado Cogstate (celdas en blanco) y pacientes mayores."
#  - Fill in the blank Cogstate result cells (HR2:JI2, minus a few that
#    legitimately stay blank) with their real values.
#  - Correct the patient's birth date (older patient: 1931, not 1941).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix patient birth date (E2): 25/02/1941 -> 25/02/1931 --------------
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "25/02/1931"

# --- Fill the previously blank Cogstate cells ---------------------------
$cogstateValues = @{
    "HR2" = "23075"
    "HS2" = "0.61548"
    "HT2" = "4"
    "HU2" = "0"
    "HV2" = "12"
    "HW2" = "31857"
    "HX2" = "0.78540"
    "HY2" = "6"
    "HZ2" = "0"
    "IA2" = "12"
    "IB2" = "37185"
    "IC2" = "0.78540"
    "ID2" = "6"
    "IE2" = "1"
    "IF2" = "12"
    "IG2" = "0"
    "IH2" = "38"
    "IK2" = "0.00000"
    "IL2" = "0"
    "IM2" = "50"
    "IP2" = "0.00000"
    "IQ2" = "36"
    "IR2" = "52"
    "IS2" = "318104"
    "IT2" = "0.16336"
    "IU2" = "0.69398"
    "IV2" = "8"
    "IW2" = "36"
    "IX2" = "322266"
    "IY2" = "0.10278"
    "IZ2" = "0.44051"
    "JA2" = "20"
    "JB2" = "27"
    "JC2" = "390424"
    "JD2" = "0.18649"
    "JE2" = "0.71065"
    "JF2" = "7"
    "JG2" = "1"
    "JH2" = "0.86912"
    "JI2" = "67830"
}

# Force text storage (so numeric-looking strings like "0.00000" or
# "23075" are kept verbatim as shared-string text, matching every other
# cell on this sheet) and then drop the temporary number-format override
# so the cells end up back on the sheet's default style.
$valuesRange = $ws.Range("HR2:JI2")
$valuesRange.NumberFormat = "@"
foreach ($addr in $cogstateValues.Keys) {
    $ws.Range($addr).Value = $cogstateValues[$addr]
}
$valuesRange.ClearFormats()

$ws.Range("E2").ClearFormats()
